# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Ifrit_Profits workbook
# (currentAveragePrice / LevePrice / LeveProfit columns per job sheet).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1030.1
$ws.Range("I18").Value = 1022.3333
$ws.Range("K18").Value = 1022.3333
$ws.Range("M18").Value = -738.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3176.5122
$ws.Range("I76").Value = 3165.7568
$ws.Range("J76").Value = 3276
$ws.Range("K76").Value = 3165.7568
$ws.Range("L76").Value = 3276
$ws.Range("M76").Value = -2850.7568
$ws.Range("N76").Value = -3906

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3176.5122
$ws.Range("I79").Value = 3165.7568
$ws.Range("J79").Value = 3276
$ws.Range("K79").Value = 3165.7568
$ws.Range("L79").Value = 3276
$ws.Range("M79").Value = -2073.7568
$ws.Range("N79").Value = -5460

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5057.5293
$ws.Range("I86").Value = 2549.5
$ws.Range("J86").Value = 5829.231
$ws.Range("K86").Value = 2549.5
$ws.Range("L86").Value = 5829.231
$ws.Range("M86").Value = -1426.5
$ws.Range("N86").Value = -8075.231

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 5057.5293
$ws.Range("I89").Value = 2549.5
$ws.Range("J89").Value = 5829.231
$ws.Range("K89").Value = 12747.5
$ws.Range("L89").Value = 29146.155
$ws.Range("M89").Value = -7131.5
$ws.Range("N89").Value = -40378.155

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 183918.66
$ws.Range("I132").Value = 194496.47
$ws.Range("J132").Value = 570
$ws.Range("K132").Value = 583489.41
$ws.Range("L132").Value = 1710
$ws.Range("M132").Value = -580959.41
$ws.Range("N132").Value = -6770

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 37039200
$ws.Range("I137").Value = 1335
$ws.Range("J137").Value = 200005800
$ws.Range("K137").Value = 4005
$ws.Range("L137").Value = 600017400
$ws.Range("M137").Value = -1455
$ws.Range("N137").Value = -600022500

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1889116.6
$ws.Range("I138").Value = 2779469.8
$ws.Range("J138").Value = 3662.8235
$ws.Range("K138").Value = 8338409.399999999
$ws.Range("L138").Value = 10988.4705
$ws.Range("M138").Value = -8333269.399999999
$ws.Range("N138").Value = -21268.4705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 817832.4
$ws.Range("I2").Value = 696.13635
$ws.Range("J2").Value = 2101903.5
$ws.Range("K2").Value = 696.13635
$ws.Range("L2").Value = 2101903.5
$ws.Range("M2").Value = -583.13635
$ws.Range("N2").Value = -2102129.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6364.702
$ws.Range("I32").Value = 7197.3716
$ws.Range("K32").Value = 7197.3716
$ws.Range("M32").Value = -6910.3716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1330.2258
$ws.Range("I61").Value = 1311.963
$ws.Range("J61").Value = 1453.5
$ws.Range("K61").Value = 1311.963
$ws.Range("L61").Value = 1453.5
$ws.Range("M61").Value = -1099.963
$ws.Range("N61").Value = -1877.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4458.162
$ws.Range("I74").Value = 1001.0417
$ws.Range("J74").Value = 10840.538
$ws.Range("K74").Value = 1001.0417
$ws.Range("L74").Value = 10840.538
$ws.Range("M74").Value = -127.0417
$ws.Range("N74").Value = -12588.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4458.162
$ws.Range("I77").Value = 1001.0417
$ws.Range("J77").Value = 10840.538
$ws.Range("K77").Value = 5005.2085
$ws.Range("L77").Value = 54202.69
$ws.Range("M77").Value = -637.2084999999997
$ws.Range("N77").Value = -62938.69

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 817832.4
$ws.Range("I116").Value = 696.13635
$ws.Range("J116").Value = 2101903.5
$ws.Range("K116").Value = 696.13635
$ws.Range("L116").Value = 2101903.5
$ws.Range("M116").Value = 1597.86365
$ws.Range("N116").Value = -2106491.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 31553.795
$ws.Range("I132").Value = 2009.6538
$ws.Range("J132").Value = 127572.25
$ws.Range("K132").Value = 6028.9614
$ws.Range("L132").Value = 382716.75
$ws.Range("M132").Value = -3498.9614
$ws.Range("N132").Value = -387776.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1330.2258
$ws.Range("I136").Value = 1311.963
$ws.Range("J136").Value = 1453.5
$ws.Range("K136").Value = 3935.889
$ws.Range("L136").Value = 4360.5
$ws.Range("M136").Value = -1385.889
$ws.Range("N136").Value = -9460.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 817832.4
$ws.Range("I3").Value = 696.13635
$ws.Range("J3").Value = 2101903.5
$ws.Range("K3").Value = 696.13635
$ws.Range("L3").Value = 2101903.5
$ws.Range("M3").Value = -582.13635
$ws.Range("N3").Value = -2102131.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 21042.428
$ws.Range("J100").Value = 21042.428
$ws.Range("L100").Value = 21042.428
$ws.Range("N100").Value = -23206.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 39350.516
$ws.Range("I134").Value = 42080.184
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 126240.552
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -123705.552
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1417.6578
$ws.Range("I31").Value = 1112.6666
$ws.Range("J31").Value = 1794.4117
$ws.Range("K31").Value = 1112.6666
$ws.Range("L31").Value = 1794.4117
$ws.Range("M31").Value = -817.6666
$ws.Range("N31").Value = -2384.4117

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1417.6578
$ws.Range("I34").Value = 1112.6666
$ws.Range("J34").Value = 1794.4117
$ws.Range("K34").Value = 1112.6666
$ws.Range("L34").Value = 1794.4117
$ws.Range("M34").Value = -910.6666
$ws.Range("N34").Value = -2198.4117

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 901.44446
$ws.Range("I122").Value = 862.1429000000001
$ws.Range("J122").Value = 1039
$ws.Range("K122").Value = 2586.4287
$ws.Range("L122").Value = 3117
$ws.Range("M122").Value = -136.4287000000004
$ws.Range("N122").Value = -8017

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1308.2683
$ws.Range("I132").Value = 1082.6364
$ws.Range("J132").Value = 2239
$ws.Range("K132").Value = 3247.9092
$ws.Range("L132").Value = 6717
$ws.Range("M132").Value = -717.9092000000001
$ws.Range("N132").Value = -11777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 1000
$ws.Range("I111").Value = 1000
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3000
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 67
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 895.85
$ws.Range("J129").Value = 1126.3636
$ws.Range("L129").Value = 3379.0908
$ws.Range("N129").Value = -13379.0908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1145.6052
$ws.Range("I102").Value = 1044.6
$ws.Range("J102").Value = 1339.8462
$ws.Range("K102").Value = 1044.6
$ws.Range("L102").Value = 1339.8462
$ws.Range("M102").Value = 577.4000000000001
$ws.Range("N102").Value = -4583.8462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1631.2222
$ws.Range("I126").Value = 1550.7059
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 4652.1177
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -2182.1177
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1836.2439
$ws.Range("I132").Value = 1685.4482
$ws.Range("J132").Value = 2200.6667
$ws.Range("K132").Value = 5056.3446
$ws.Range("L132").Value = 6602.000100000001
$ws.Range("M132").Value = -2526.3446
$ws.Range("N132").Value = -11662.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1087.4375
$ws.Range("I61").Value = 1008.68
$ws.Range("J61").Value = 1368.7142
$ws.Range("K61").Value = 1008.68
$ws.Range("L61").Value = 1368.7142
$ws.Range("M61").Value = -806.6799999999999
$ws.Range("N61").Value = -1772.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1087.4375
$ws.Range("I113").Value = 1008.68
$ws.Range("J113").Value = 1368.7142
$ws.Range("K113").Value = 1008.68
$ws.Range("L113").Value = 1368.7142
$ws.Range("M113").Value = 1161.32
$ws.Range("N113").Value = -5708.7142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4329.021
$ws.Range("I132").Value = 4860
$ws.Range("J132").Value = 2311.3
$ws.Range("K132").Value = 14580
$ws.Range("L132").Value = 6933.900000000001
$ws.Range("M132").Value = -12050
$ws.Range("N132").Value = -11993.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9894.091
$ws.Range("I136").Value = 13097.1875
$ws.Range("J136").Value = 1352.5
$ws.Range("K136").Value = 39291.5625
$ws.Range("L136").Value = 4057.5
$ws.Range("M136").Value = -36741.5625
$ws.Range("N136").Value = -9157.5
